$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5
$ws.Range("L2").Value = "stimuli/img_mgnmm.png"
$ws.Range("M2").Value = 79.1470588235294
$ws.Range("N2").Value = 60.38235294117647
$ws.Range("O2").Value = 69.76470588235294
$ws.Range("P2").Value = 34
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = 8
$ws.Range("S2").Value = 8
$ws.Range("T2").Value = 8
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 8
$ws.Range("C3").Value = 5
$ws.Range("H3").Value = "kitchens"
$ws.Range("I3").Value = "target"
$ws.Range("K3").Value = "j"
$ws.Range("L3").Value = "stimuli/img_esb4r.png"
$ws.Range("M3").Value = 60.73529411764706
$ws.Range("N3").Value = 38.58823529411764
$ws.Range("O3").Value = 49.66176470588235
$ws.Range("P3").Value = 34
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 3
$ws.Range("V3").Value = 3
$ws.Range("C4").Value = 5
$ws.Range("H4").Value = "kitchens"
$ws.Range("I4").Value = "target"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_1ao2d.png"
$ws.Range("M4").Value = 38.77777777777778
$ws.Range("N4").Value = 18.75
$ws.Range("O4").Value = 28.76388888888889
$ws.Range("P4").Value = 36
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 1
$ws.Range("C5").Value = 5
$ws.Range("H5").Value = "bedrooms"
$ws.Range("I5").Value = "distractor"
$ws.Range("K5").Value = "f"
$ws.Range("L5").Value = "stimuli/img_n9xll.png"
$ws.Range("M5").Value = 77.14285714285714
$ws.Range("N5").Value = 59.21428571428572
$ws.Range("O5").Value = 68.17857142857143
$ws.Range("P5").Value = 42
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = 7
$ws.Range("S5").Value = 7
$ws.Range("T5").Value = 7
$ws.Range("U5").Value = 7
$ws.Range("V5").Value = 7
$ws.Range("C6").Value = 5
$ws.Range("H6").Value = "bedrooms"
$ws.Range("L6").Value = "stimuli/img_dkqas.png"
$ws.Range("M6").Value = 78.57894736842105
$ws.Range("N6").Value = 57.71052631578947
$ws.Range("O6").Value = 68.14473684210526
$ws.Range("P6").Value = 38
$ws.Range("Q6").Value = 7
$ws.Range("R6").Value = 7
$ws.Range("S6").Value = 7
$ws.Range("T6").Value = 7
$ws.Range("U6").Value = 7
$ws.Range("V6").Value = 7
$ws.Range("C7").Value = 5
$ws.Range("H7").Value = "living_rooms"
$ws.Range("I7").Value = "distractor"
$ws.Range("K7").Value = "f"
$ws.Range("L7").Value = "stimuli/img_pna7l.png"
$ws.Range("M7").Value = 85.53333333333333
$ws.Range("N7").Value = 67.97777777777777
$ws.Range("O7").Value = 76.75555555555556
$ws.Range("P7").Value = 45
$ws.Range("C8").Value = 5
$ws.Range("H8").Value = "kitchens"
$ws.Range("I8").Value = "target"
$ws.Range("K8").Value = "j"
$ws.Range("L8").Value = "stimuli/img_wppku.png"
$ws.Range("M8").Value = 75.02941176470588
$ws.Range("N8").Value = 53.05882352941177
$ws.Range("O8").Value = 64.04411764705883
$ws.Range("P8").Value = 34
$ws.Range("Q8").Value = 6
$ws.Range("R8").Value = 6
$ws.Range("S8").Value = 6
$ws.Range("T8").Value = 6
$ws.Range("U8").Value = 6
$ws.Range("V8").Value = 6
$ws.Range("C9").Value = 5
$ws.Range("L9").Value = "stimuli/img_7ucnr.png"
$ws.Range("M9").Value = 70.39393939393939
$ws.Range("N9").Value = 47.90909090909091
$ws.Range("O9").Value = 59.15151515151515
$ws.Range("P9").Value = 33
$ws.Range("Q9").Value = 5
$ws.Range("R9").Value = 5
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 5
$ws.Range("U9").Value = 5
$ws.Range("V9").Value = 5
$ws.Range("C10").Value = 5
$ws.Range("L10").Value = "stimuli/img_amsgw.png"
$ws.Range("M10").Value = 86.08510638297872
$ws.Range("N10").Value = 65.95744680851064
$ws.Range("O10").Value = 76.02127659574468
$ws.Range("P10").Value = 47
$ws.Range("Q10").Value = 9
$ws.Range("R10").Value = 9
$ws.Range("S10").Value = 9
$ws.Range("T10").Value = 8
$ws.Range("U10").Value = 9
$ws.Range("V10").Value = 8
$ws.Range("C11").Value = 5
$ws.Range("H11").Value = "living_rooms"
$ws.Range("L11").Value = "stimuli/img_ac0ey.png"
$ws.Range("M11").Value = 86.62222222222222
$ws.Range("N11").Value = 70.02222222222223
$ws.Range("O11").Value = 78.32222222222222
$ws.Range("P11").Value = 45
$ws.Range("Q11").Value = 9
$ws.Range("R11").Value = 9
$ws.Range("S11").Value = 9
$ws.Range("T11").Value = 9
$ws.Range("U11").Value = 9
$ws.Range("V11").Value = 9
$ws.Range("C12").Value = 5
$ws.Range("H12").Value = "bedrooms"
$ws.Range("L12").Value = "stimuli/img_5m6x4.png"
$ws.Range("M12").Value = 80.23076923076923
$ws.Range("N12").Value = 58.41025641025641
$ws.Range("O12").Value = 69.32051282051282
$ws.Range("P12").Value = 39
$ws.Range("Q12").Value = 7
$ws.Range("R12").Value = 7
$ws.Range("S12").Value = 7
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 7
$ws.Range("V12").Value = 7
$ws.Range("C13").Value = 5
$ws.Range("H13").Value = "living_rooms"
$ws.Range("I13").Value = "distractor"
$ws.Range("K13").Value = "f"
$ws.Range("L13").Value = "stimuli/img_95hiv.png"
$ws.Range("M13").Value = 84.04545454545455
$ws.Range("N13").Value = 67.31818181818181
$ws.Range("O13").Value = 75.68181818181819
$ws.Range("P13").Value = 44
$ws.Range("Q13").Value = 9
$ws.Range("R13").Value = 9
$ws.Range("S13").Value = 9
$ws.Range("T13").Value = 8
$ws.Range("U13").Value = 8
$ws.Range("V13").Value = 8
$ws.Range("C14").Value = 5
$ws.Range("H14").Value = "kitchens"
$ws.Range("I14").Value = "target"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_ikk62.png"
$ws.Range("M14").Value = 37.48780487804878
$ws.Range("N14").Value = 21.07317073170732
$ws.Range("O14").Value = 29.28048780487805
$ws.Range("P14").Value = 41
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 1
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 1
$ws.Range("U14").Value = 1
$ws.Range("V14").Value = 1
$ws.Range("C15").Value = 5
$ws.Range("H15").Value = "bedrooms"
$ws.Range("I15").Value = "distractor"
$ws.Range("K15").Value = "f"
$ws.Range("L15").Value = "stimuli/img_u1rxv.png"
$ws.Range("M15").Value = 75.63636363636364
$ws.Range("N15").Value = 54.27272727272727
$ws.Range("O15").Value = 64.95454545454545
$ws.Range("P15").Value = 44
$ws.Range("Q15").Value = 6
$ws.Range("R15").Value = 6
$ws.Range("S15").Value = 6
$ws.Range("T15").Value = 6
$ws.Range("U15").Value = 6
$ws.Range("V15").Value = 6
$ws.Range("C16").Value = 5
$ws.Range("H16").Value = "living_rooms"
$ws.Range("I16").Value = "distractor"
$ws.Range("K16").Value = "f"
$ws.Range("L16").Value = "stimuli/img_f63yi.png"
$ws.Range("M16").Value = 85.275
$ws.Range("N16").Value = 68.475
$ws.Range("O16").Value = 76.875
$ws.Range("P16").Value = 40
$ws.Range("Q16").Value = 9
$ws.Range("R16").Value = 9
$ws.Range("S16").Value = 9
$ws.Range("T16").Value = 9
$ws.Range("U16").Value = 8
$ws.Range("V16").Value = 9
$ws.Range("C17").Value = 5
$ws.Range("L17").Value = "stimuli/img_6zz63.png"
$ws.Range("M17").Value = 87.66666666666667
$ws.Range("N17").Value = 70.6
$ws.Range("O17").Value = 79.13333333333333
$ws.Range("P17").Value = 45
$ws.Range("Q17").Value = 9
$ws.Range("R17").Value = 10
$ws.Range("S17").Value = 10
$ws.Range("T17").Value = 9
$ws.Range("U17").Value = 9
$ws.Range("V17").Value = 9
$ws.Range("C18").Value = 5
$ws.Range("L18").Value = "stimuli/img_1zhz6.png"
$ws.Range("M18").Value = 49.02272727272727
$ws.Range("N18").Value = 32.77272727272727
$ws.Range("O18").Value = 40.89772727272727
$ws.Range("P18").Value = 44
$ws.Range("Q18").Value = 3
$ws.Range("R18").Value = 3
$ws.Range("S18").Value = 3
$ws.Range("T18").Value = 3
$ws.Range("U18").Value = 3
$ws.Range("V18").Value = 3
$ws.Range("C19").Value = 5
$ws.Range("L19").Value = "stimuli/img_mawe6.png"
$ws.Range("M19").Value = 83.48387096774194
$ws.Range("N19").Value = 65.54838709677419
$ws.Range("O19").Value = 74.51612903225806
$ws.Range("P19").Value = 31
$ws.Range("Q19").Value = 9
$ws.Range("R19").Value = 9
$ws.Range("S19").Value = 9
$ws.Range("T19").Value = 9
$ws.Range("U19").Value = 9
$ws.Range("V19").Value = 9
$ws.Range("C20").Value = 5
$ws.Range("H20").Value = "bedrooms"
$ws.Range("I20").Value = "distractor"
$ws.Range("K20").Value = "f"
$ws.Range("L20").Value = "stimuli/img_d9ogj.png"
$ws.Range("M20").Value = 76.86842105263158
$ws.Range("N20").Value = 53.5
$ws.Range("O20").Value = 65.18421052631578
$ws.Range("P20").Value = 38
$ws.Range("Q20").Value = 6
$ws.Range("R20").Value = 6
$ws.Range("S20").Value = 6
$ws.Range("T20").Value = 6
$ws.Range("U20").Value = 6
$ws.Range("V20").Value = 6
$ws.Range("C21").Value = 5
$ws.Range("H21").Value = "kitchens"
$ws.Range("I21").Value = "target"
$ws.Range("K21").Value = "j"
$ws.Range("L21").Value = "stimuli/img_qmgwq.png"
$ws.Range("M21").Value = 84.58333333333333
$ws.Range("N21").Value = 64.44444444444444
$ws.Range("O21").Value = 74.51388888888889
$ws.Range("P21").Value = 36
$ws.Range("Q21").Value = 9
$ws.Range("R21").Value = 9
$ws.Range("S21").Value = 9
$ws.Range("T21").Value = 9
$ws.Range("U21").Value = 9
$ws.Range("V21").Value = 9
$ws.Range("C22").Value = 5
$ws.Range("H22").Value = "kitchens"
$ws.Range("I22").Value = "target"
$ws.Range("K22").Value = "j"
$ws.Range("L22").Value = "stimuli/img_kwxq1.png"
$ws.Range("M22").Value = 68.53125
$ws.Range("N22").Value = 44.09375
$ws.Range("O22").Value = 56.3125
$ws.Range("P22").Value = 32
$ws.Range("U22").Value = 4
$ws.Range("C23").Value = 5
$ws.Range("H23").Value = "bedrooms"
$ws.Range("L23").Value = "stimuli/img_8dacu.png"
$ws.Range("M23").Value = 76.38461538461539
$ws.Range("N23").Value = 53.64102564102564
$ws.Range("O23").Value = 65.01282051282051
$ws.Range("P23").Value = 39
$ws.Range("Q23").Value = 6
$ws.Range("R23").Value = 6
$ws.Range("S23").Value = 6
$ws.Range("T23").Value = 6
$ws.Range("U23").Value = 6
$ws.Range("V23").Value = 6
$ws.Range("C24").Value = 5
$ws.Range("H24").Value = "kitchens"
$ws.Range("I24").Value = "target"
$ws.Range("K24").Value = "j"
$ws.Range("L24").Value = "stimuli/img_zi8qc.png"
$ws.Range("M24").Value = 77.14285714285714
$ws.Range("N24").Value = 57.02857142857143
$ws.Range("O24").Value = 67.08571428571429
$ws.Range("P24").Value = 35
$ws.Range("Q24").Value = 7
$ws.Range("R24").Value = 7
$ws.Range("S24").Value = 7
$ws.Range("T24").Value = 7
$ws.Range("U24").Value = 7
$ws.Range("V24").Value = 7
$ws.Range("C25").Value = 5
$ws.Range("L25").Value = "stimuli/img_89rmb.png"
$ws.Range("M25").Value = 55.18518518518518
$ws.Range("N25").Value = 29.25925925925926
$ws.Range("O25").Value = 42.22222222222222
$ws.Range("P25").Value = 27
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 2
$ws.Range("S25").Value = 2
$ws.Range("T25").Value = 2
$ws.Range("U25").Value = 2
$ws.Range("V25").Value = 2
$ws.Range("C26").Value = 5
$ws.Range("H26").Value = "kitchens"
$ws.Range("I26").Value = "target"
$ws.Range("K26").Value = "j"
$ws.Range("L26").Value = "stimuli/img_xdhz2.png"
$ws.Range("M26").Value = 63.3
$ws.Range("N26").Value = 37.25
$ws.Range("O26").Value = 50.275
$ws.Range("P26").Value = 40
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = 3
$ws.Range("S26").Value = 3
$ws.Range("T26").Value = 3
$ws.Range("U26").Value = 3
$ws.Range("V26").Value = 3
$ws.Range("C27").Value = 5
$ws.Range("H27").Value = "bedrooms"
$ws.Range("I27").Value = "distractor"
$ws.Range("K27").Value = "f"
$ws.Range("L27").Value = "stimuli/img_d26ik.png"
$ws.Range("M27").Value = 77.73809523809524
$ws.Range("N27").Value = 60.66666666666666
$ws.Range("O27").Value = 69.20238095238095
$ws.Range("P27").Value = 42
$ws.Range("Q27").Value = 7
$ws.Range("R27").Value = 7
$ws.Range("S27").Value = 7
$ws.Range("T27").Value = 7
$ws.Range("U27").Value = 7
$ws.Range("V27").Value = 7
$ws.Range("C28").Value = 5
$ws.Range("L28").Value = "stimuli/img_7w5tw.png"
$ws.Range("M28").Value = 53.2258064516129
$ws.Range("N28").Value = 28.90322580645161
$ws.Range("O28").Value = 41.06451612903226
$ws.Range("P28").Value = 31
$ws.Range("Q28").Value = 2
$ws.Range("R28").Value = 2
$ws.Range("S28").Value = 2
$ws.Range("T28").Value = 2
$ws.Range("U28").Value = 2
$ws.Range("V28").Value = 2
$ws.Range("C29").Value = 5
$ws.Range("L29").Value = "stimuli/img_ewrjk.png"
$ws.Range("M29").Value = 73.0909090909091
$ws.Range("N29").Value = 53.39393939393939
$ws.Range("O29").Value = 63.24242424242424
$ws.Range("P29").Value = 33
$ws.Range("Q29").Value = 6
$ws.Range("R29").Value = 6
$ws.Range("S29").Value = 6
$ws.Range("T29").Value = 6
$ws.Range("U29").Value = 6
$ws.Range("V29").Value = 6
$ws.Range("C30").Value = 5
$ws.Range("L30").Value = "stimuli/img_lszzj.png"
$ws.Range("M30").Value = 64.70588235294117
$ws.Range("N30").Value = 45.58823529411764
$ws.Range("O30").Value = 55.14705882352941
$ws.Range("P30").Value = 34
$ws.Range("Q30").Value = 4
$ws.Range("R30").Value = 4
$ws.Range("S30").Value = 4
$ws.Range("T30").Value = 4
$ws.Range("U30").Value = 4
$ws.Range("V30").Value = 4
$ws.Range("C31").Value = 5
$ws.Range("L31").Value = "stimuli/img_mjxmq.png"
$ws.Range("M31").Value = 77.07692307692308
$ws.Range("N31").Value = 58.15384615384615
$ws.Range("O31").Value = 67.61538461538461
$ws.Range("P31").Value = 39
$ws.Range("C32").Value = 5
$ws.Range("H32").Value = "kitchens"
$ws.Range("I32").Value = "target"
$ws.Range("K32").Value = "j"
$ws.Range("L32").Value = "stimuli/img_vbrb7.png"
$ws.Range("M32").Value = 85.5625
$ws.Range("N32").Value = 71.46875
$ws.Range("O32").Value = 78.515625
$ws.Range("P32").Value = 32
$ws.Range("Q32").Value = 10
$ws.Range("R32").Value = 10
$ws.Range("S32").Value = 10
$ws.Range("T32").Value = 10
$ws.Range("U32").Value = 10
$ws.Range("V32").Value = 10
$ws.Range("C33").Value = 5
$ws.Range("H33").Value = "kitchens"
$ws.Range("I33").Value = "target"
$ws.Range("K33").Value = "j"
$ws.Range("L33").Value = "stimuli/img_r2lxk.png"
$ws.Range("M33").Value = 89.24242424242425
$ws.Range("N33").Value = 67.6969696969697
$ws.Range("O33").Value = 78.46969696969697
$ws.Range("P33").Value = 33
$ws.Range("Q33").Value = 10
$ws.Range("R33").Value = 10
$ws.Range("S33").Value = 10
$ws.Range("T33").Value = 10
$ws.Range("U33").Value = 10
$ws.Range("V33").Value = 10
$ws.Range("C34").Value = 5
$ws.Range("H34").Value = "living_rooms"
$ws.Range("L34").Value = "stimuli/img_yosqb.png"
$ws.Range("M34").Value = 50.88372093023256
$ws.Range("N34").Value = 30.11627906976744
$ws.Range("O34").Value = 40.5
$ws.Range("P34").Value = 43
$ws.Range("Q34").Value = 3
$ws.Range("R34").Value = 3
$ws.Range("S34").Value = 3
$ws.Range("T34").Value = 3
$ws.Range("U34").Value = 3
$ws.Range("V34").Value = 3
$ws.Range("C35").Value = 5
$ws.Range("H35").Value = "living_rooms"
$ws.Range("I35").Value = "distractor"
$ws.Range("K35").Value = "f"
$ws.Range("L35").Value = "stimuli/img_0jzz7.png"
$ws.Range("M35").Value = 84.85106382978724
$ws.Range("N35").Value = 68.87234042553192
$ws.Range("O35").Value = 76.86170212765958
$ws.Range("P35").Value = 47
$ws.Range("Q35").Value = 9
$ws.Range("R35").Value = 9
$ws.Range("S35").Value = 9
$ws.Range("T35").Value = 9
$ws.Range("U35").Value = 8
$ws.Range("V35").Value = 9
$ws.Range("C36").Value = 5
$ws.Range("H36").Value = "bedrooms"
$ws.Range("I36").Value = "distractor"
$ws.Range("K36").Value = "f"
$ws.Range("L36").Value = "stimuli/img_ybbmx.png"
$ws.Range("M36").Value = 55.24324324324324
$ws.Range("N36").Value = 36.75675675675676
$ws.Range("O36").Value = 46
$ws.Range("P36").Value = 37
$ws.Range("Q36").Value = 3
$ws.Range("R36").Value = 3
$ws.Range("S36").Value = 3
$ws.Range("T36").Value = 3
$ws.Range("U36").Value = 3
$ws.Range("V36").Value = 3
$ws.Range("C37").Value = 5
$ws.Range("L37").Value = "stimuli/img_fea1z.png"
$ws.Range("M37").Value = 79.45945945945945
$ws.Range("N37").Value = 56.24324324324324
$ws.Range("O37").Value = 67.85135135135135
$ws.Range("P37").Value = 37
$ws.Range("Q37").Value = 7
$ws.Range("R37").Value = 7
$ws.Range("S37").Value = 7
$ws.Range("T37").Value = 7
$ws.Range("U37").Value = 7
$ws.Range("V37").Value = 7
$ws.Range("C38").Value = 5
$ws.Range("H38").Value = "bedrooms"
$ws.Range("I38").Value = "distractor"
$ws.Range("K38").Value = "f"
$ws.Range("L38").Value = "stimuli/img_iqmdm.png"
$ws.Range("M38").Value = 79.38888888888889
$ws.Range("N38").Value = 58.36111111111111
$ws.Range("O38").Value = 68.875
$ws.Range("P38").Value = 36
$ws.Range("Q38").Value = 7
$ws.Range("R38").Value = 7
$ws.Range("S38").Value = 7
$ws.Range("T38").Value = 7
$ws.Range("U38").Value = 7
$ws.Range("V38").Value = 7
$ws.Range("C39").Value = 5
$ws.Range("L39").Value = "stimuli/img_7ed9m.png"
$ws.Range("M39").Value = 80.71875
$ws.Range("N39").Value = 58.65625
$ws.Range("O39").Value = 69.6875
$ws.Range("P39").Value = 32
$ws.Range("Q39").Value = 8
$ws.Range("R39").Value = 8
$ws.Range("S39").Value = 8
$ws.Range("T39").Value = 8
$ws.Range("U39").Value = 8
$ws.Range("V39").Value = 8
$ws.Range("C40").Value = 5
$ws.Range("L40").Value = "stimuli/img_wgkqa.png"
$ws.Range("M40").Value = 87.25581395348837
$ws.Range("N40").Value = 71.13953488372093
$ws.Range("O40").Value = 79.19767441860465
$ws.Range("P40").Value = 43
$ws.Range("Q40").Value = 10
$ws.Range("R40").Value = 10
$ws.Range("S40").Value = 10
$ws.Range("T40").Value = 9
$ws.Range("U40").Value = 9
$ws.Range("V40").Value = 9
$ws.Range("C41").Value = 5
$ws.Range("L41").Value = "stimuli/img_z293c.png"
$ws.Range("M41").Value = 71.26470588235294
$ws.Range("N41").Value = 46.88235294117647
$ws.Range("O41").Value = 59.07352941176471
$ws.Range("P41").Value = 34
$ws.Range("Q41").Value = 5
$ws.Range("R41").Value = 5
$ws.Range("S41").Value = 5
$ws.Range("T41").Value = 5
$ws.Range("U41").Value = 5
$ws.Range("V41").Value = 5

Write-Host "Applied changes:" 528